$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply each cell update from the diff. Numeric-looking text values
# (e.g. "131.00", "0.0000171") must stay as TEXT, matching the source
# inline strings -- force the "@" text format before assignment, then
# drop back to the default "Normal" style so no stray formatting is left
# behind (matches cells that carry no explicit style in the original file).

$ws.Range("D2").Value = "62.494.85"
$ws.Range("E2").Value = "  +2.15%  "
$ws.Range("D3").Value = "3.432.27"
$ws.Range("E3").Value = "  +3.02%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "407.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.66%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.600"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.31%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.698"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.144"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +22.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.34"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.20%  "
$ws.Range("E12").Value = "  +0.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.51"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.00%  "
$ws.Range("D15").Value = "3.416.14"
$ws.Range("E15").Value = "  +2.33%  "
$ws.Range("D16").Value = "62.501.79"
$ws.Range("E16").Value = "  +2.11%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000171"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +35.71%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.57%  "
$ws.Range("B19").Value = "Polygon"
$ws.Range("C19").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.02"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "84.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "315.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.79%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.26%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "29.85"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.19"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.28%  "
$ws.Range("E28").Value = "  +5.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +9.07%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.174"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.64%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "44.37"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.42%  "
$ws.Range("E32").Value = "  +2.57%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.45"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0487"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "51.78"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.15%  "
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.98"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.55%  "
$ws.Range("E39").Value = "  +16.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.33"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "143.66"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.49%  "
$ws.Range("E42").Value = "  +3.61%  "
$ws.Range("E43").Value = "  +2.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.02"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.92%  "
$ws.Range("E45").Value = "  +2.53%  "
$ws.Range("E46").Value = "  +0.26%  "
$ws.Range("E47").Value = "  +1.78%  "
$ws.Range("D48").Value = "2.114.91"
$ws.Range("E48").Value = "  +0.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +10.14%  "
$ws.Range("E50").Value = "  -0.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.09"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +32.00%  "
